$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New kingdom passive skill rows (22-27), continuing the existing table.
$newRows = @(
    @{ Row=22; A=21; B="Master Farmer";      C="As you level this skill you will get 10% extra population on the hourly update. At level 10, you will gain 100% more population."; D=10; E=1; F=0.1; K=15; L=8;  M=5; N=1; O=1 },
    @{ Row=23; A=22; B="Master Stone Mason"; C="Gain an additonal 10% stone per level for an additional 100% when the hourly update hits."; D=10; E=1; F=0.1; K=16; L=15; M=4; N=1; O=1 },
    @{ Row=24; A=23; B="Master Wood Worker"; C="Gain an extra 10% wood for a max of 100% on the hourly update."; D=10; E=1; F=0.1; K=17; L=2;  M=3; N=1; O=1 },
    @{ Row=25; A=24; B="Master of Iron";     C="As you level this skill you will gain an additional 10% per level of extra iron for a max of 100% at max level when the hourly update hits."; D=10; E=2; F=0.1; K=18; L=6;  M=5; N=1; O=1 },
    @{ Row=26; A=25; B="Master of Steel";    C="You will gain an additional 10% steel per level for a max of 100% when you smelt steel."; D=10; E=3; F=0.1; K=20; L=11; M=5; N=1; O=1 },
    @{ Row=27; A=26; B="Master Potter";      C="As you level this skill you will gain 10% more clay per hour on the hourly update for a max of 100% at max level."; D=10; E=1; F=0.1; K=19; L=2;  M=3; N=1; O=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
}
